$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.440.05"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "2.248.41"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.33"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.87"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.01"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.01"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0812"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.22"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").Value = "  +1.21%  "
$ws.Range("D14").Value = "2.357.40"
$ws.Range("E14").Value = "  +4.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.843"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.72"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "44.133.36"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.38"
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.41"
$ws.Range("E20").Value = "  +2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "65.91"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("B22").Value = "PancakeSwap"
$ws.Range("C22").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.03"
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.65"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("E24").Value = "  +4.44%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.42"
$ws.Range("E26").Value = "  +6.47%  "
$ws.Range("E27").Value = "  +4.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.89"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.98"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.11"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.26"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0802"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.12"
$ws.Range("E34").Value = "  -7.60%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.120"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  +3.29%  "
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("E38").Value = "  +5.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.76"
$ws.Range("E39").Value = "  -1.37%  "
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("E41").Value = "  +1.05%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "1.741.48"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("E44").Value = "  +3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "80.77"
$ws.Range("E45").Value = "  -5.36%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.18"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.97"
$ws.Range("E47").Value = "  -2.64%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "71.00"
$ws.Range("E48").Value = "  +4.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "56.29"
$ws.Range("E49").Value = "  +3.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.60"
$ws.Range("E50").Value = "  +6.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.13"
$ws.Range("E51").Value = "  -0.30%  "
